# Adding 2080 to the lcoe table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# The existing last row of the table (row 9) is duplicated down to a new
# row 10, carrying forward its original values (including Year = 0), and
# then row 9's Year value is updated to the newly added 2080 entry.
$ws.Range("D9:L9").Copy($ws.Range("D10")) | Out-Null
$ws.Range("F9").Value = 2080

$ws.Range("I13").Select() | Out-Null
